$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.0003704876871779561

$ws.Range("A3").Value = 0.00011398269271012396
$ws.Range("C3").Value = 44.230770111083984
$ws.Range("D3").Value = 22.11538314819336

$ws.Range("A4").Value = 0.00010937692422885448
$ws.Range("C4").Value = 92.69230651855469
$ws.Range("D4").Value = 46.346153259277344

$ws.Range("A5").Value = 0.00003586153979995288
$ws.Range("C5").Value = 48.46154022216797
$ws.Range("D5").Value = 24.23076820373535

$ws.Range("A6").Value = 0.00003091499820584431
$ws.Range("C6").Value = 44.03845977783203
$ws.Range("D6").Value = 22.104520797729492

$ws.Range("A7").Value = 0.000026775000151246786
$ws.Range("C7").Value = 48.46154022216797
$ws.Range("D7").Value = 24.23076820373535

$ws.Range("A8").Value = 0.000021275000108289532
$ws.Range("C8").Value = 44.230770111083984
$ws.Range("D8").Value = 22.11538314819336

$ws.Range("A9").Value = 0.000010637308150762692
$ws.Range("C9").Value = 48.46154022216797
$ws.Range("D9").Value = 22.11538314819336

$ws.Range("A10").Value = 0.000008890384378901217
$ws.Range("C10").Value = 44.230770111083984
$ws.Range("D10").Value = 22.11538314819336

$ws.Range("A12").Value = 0.000005935769422649173
$ws.Range("C12").Value = 48.653846740722656
$ws.Range("D12").Value = 24.352018356323242

